# Adding Screenshot methods. Code improvisation.
#
# The "checkoutWithExpensiveItem" smoke-test row (row 2) is updated to point
# at a new product/price fixture, and a new "currency" column is introduced
# so the test data records which currency format the price string uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SmokeTests")

# New header in column G: "currency"
$ws.Range("G1").Value = "currency"

# productName (E2): swap the old Samsung TV listing for the new Sony one
$ws.Range("E2").Value = "SONY 65 INCHES 65X7000E 4k UHD HDR SMART LED TV + ONE YEAR DEALER'S WARRANTY"

# productPrice (F2): new price, now formatted with an "Rs." prefix instead of the rupee sign
$ws.Range("F2").Value = "Rs.153,799"

# Reflect the new selection left by the edit (was J11, now G2 - the new currency cell)
$ws.Range("G2").Select()
